$d = $word.ActiveDocument
$bullet = [char]8226

# --- 1. Condense the CORE COMPETENCIES section -------------------------------
# The section currently holds three long paragraphs (Product Management &
# Strategy / Technical Product Development / Platform & Infrastructure), each
# packed with detailed sub-bullets. Collapse them into a single short summary
# paragraph and remove the other two paragraphs entirely.
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("CORE COMPETENCIES")) {
        $headingPara = $i
        break
    }
}

$firstPara = $d.Paragraphs.Item($headingPara + 1)

# Replace the first paragraph's text with the condensed summary line.
$firstPara.Range.Text = "Product Management & Strategy $bullet Technical Product Development $bullet Platform & Infrastructure"

# Delete the second and third paragraphs (including their paragraph marks).
# Re-fetch indices after each delete since the collection shifts.
$d.Paragraphs.Item($headingPara + 2).Range.Delete()
$d.Paragraphs.Item($headingPara + 2).Range.Delete()

# --- 2. Append a new "TECHNICAL SKILLS" section at the end of the document ---
$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Text = "TECHNICAL SKILLS"
$p.Style = "Heading2"

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics"

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration"

$end = $d.Content.End
$r = $d.Range($end, $end)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Style = "Normal"
$p.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training"

Write-Output "done"
